$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.868.18"
$ws.Cells.Item(2, 5).Value = "  -6.45%  "
$ws.Cells.Item(3, 4).Value = "3.699.37"
$ws.Cells.Item(3, 5).Value = "  -5.88%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "578.81"
$ws.Cells.Item(5, 5).Value = "  -3.33%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "174.94"
$ws.Cells.Item(6, 5).Value = "  +2.71%  "
$ws.Cells.Item(7, 4).Value = "3.697.99"
$ws.Cells.Item(7, 5).Value = "  -5.78%  "
$ws.Cells.Item(8, 5).Value = "  -8.42%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.997"
$ws.Cells.Item(9, 5).Value = "  -0.27%  "
$ws.Cells.Item(10, 5).Value = "  -9.76%  "
$ws.Cells.Item(11, 5).Value = "  -13.39%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "51.41"
$ws.Cells.Item(12, 5).Value = "  -7.86%  "
$ws.Cells.Item(13, 5).Value = "  -12.52%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "10.39"
$ws.Cells.Item(14, 5).Value = "  -9.13%  "
$ws.Cells.Item(15, 4).Value = "4.296.83"
$ws.Cells.Item(15, 5).Value = "  -5.58%  "
$ws.Cells.Item(16, 4).Value = "3.694.22"
$ws.Cells.Item(16, 5).Value = "  -6.35%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "19.30"
$ws.Cells.Item(17, 5).Value = "  -9.34%  "
$ws.Cells.Item(18, 5).Value = "  -3.09%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.80"
$ws.Cells.Item(19, 5).Value = "  -9.26%  "
$ws.Cells.Item(20, 5).Value = "  -8.94%  "
$ws.Cells.Item(21, 4).Value = "67.670.63"
$ws.Cells.Item(21, 5).Value = "  -6.66%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "404.30"
$ws.Cells.Item(22, 5).Value = "  -9.87%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.46"
$ws.Cells.Item(23, 5).Value = "  -6.39%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "87.83"
$ws.Cells.Item(24, 5).Value = "  -7.76%  "
$ws.Cells.Item(25, 5).Value = "  -8.80%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "12.68"
$ws.Cells.Item(26, 5).Value = "  -9.47%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.71"
$ws.Cells.Item(27, 5).Value = "  -3.48%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.03"
$ws.Cells.Item(28, 5).Value = "  +1.30%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "3.77"
$ws.Cells.Item(29, 5).Value = "  -11.33%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.42"
$ws.Cells.Item(30, 5).Value = "  -8.75%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "32.42"
$ws.Cells.Item(31, 5).Value = "  -9.10%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "7.40"
$ws.Cells.Item(32, 5).Value = "  -5.59%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "12.40"
$ws.Cells.Item(33, 5).Value = "  -9.78%  "
$ws.Cells.Item(34, 2).Value = "Bittensor"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "611.47"
$ws.Cells.Item(34, 5).Value = "  -2.64%  "
$ws.Cells.Item(35, 2).Value = "OKB"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "64.82"
$ws.Cells.Item(35, 5).Value = "  -5.99%  "
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.115"
$ws.Cells.Item(36, 5).Value = "  -9.24%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "42.83"
$ws.Cells.Item(37, 5).Value = "  -15.42%  "
$ws.Cells.Item(38, 4).Value = "0.0₃0882"
$ws.Cells.Item(38, 5).Value = "  -9.57%  "
$ws.Cells.Item(39, 5).Value = "  +0.21%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.394"
$ws.Cells.Item(40, 5).Value = "  -7.34%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.00"
$ws.Cells.Item(41, 5).Value = "  -0.10%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.135"
$ws.Cells.Item(42, 5).Value = "  -6.72%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.75"
$ws.Cells.Item(43, 5).Value = "  +6.24%  "
$ws.Cells.Item(44, 5).Value = "  -10.83%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0432"
$ws.Cells.Item(45, 5).Value = "  -9.03%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.86"
$ws.Cells.Item(46, 5).Value = "  -11.12%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.19"
$ws.Cells.Item(47, 5).Value = "  -12.51%  "
$ws.Cells.Item(48, 4).Value = "2.798.75"
$ws.Cells.Item(48, 5).Value = "  -1.12%  "
$ws.Cells.Item(49, 5).Value = "  -9.30%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.69"
$ws.Cells.Item(50, 5).Value = "  -5.53%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "3.03"
$ws.Cells.Item(51, 5).Value = "  -10.05%  "
